$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.623.37'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '1.840.59'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '311.75'
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").Value = '0.4259'
$ws.Range("E7").Value = '  +0.85%  '

$ws.Range("D8").Value = '0.3604'
$ws.Range("E8").Value = '  -0.90%  '

$ws.Range("D9").Value = '0.07294'
$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("D10").Value = '0.8663'
$ws.Range("E10").Value = '  -2.56%  '

$ws.Range("D11").Value = '20.56'
$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("D12").Value = '1.829.27'
$ws.Range("E12").Value = '  -1.78%  '

$ws.Range("D13").Value = '5.319'
$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("D14").Value = '6.490'
$ws.Range("E14").Value = '  -1.14%  '

$ws.Range("D15").Value = '0.06971'
$ws.Range("E15").Value = '  +1.62%  '

$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("D17").Value = '79.26'
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").Value = '0.000008926'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.01%  '

$ws.Range("D20").Value = '15.22'
$ws.Range("E20").Value = '  -1.39%  '

$ws.Range("D21").Value = '27.620.69'
$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").Value = '4.964'
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("D23").Value = '10.32'
$ws.Range("E23").Value = '  -1.97%  '

$ws.Range("D24").Value = '2.069.17'
$ws.Range("E24").Value = '  -0.63%  '

$ws.Range("D25").Value = '1.976'
$ws.Range("E25").Value = '  +2.23%  '

$ws.Range("D26").Value = '155.35'
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").Value = '119.45'
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").Value = '5.204'
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("D30").Value = '1.861'
$ws.Range("E30").Value = '  +0.62%  '

$ws.Range("D31").Value = '0.08880'
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").Value = '0.7595'
$ws.Range("E32").Value = '  -2.38%  '

$ws.Range("D33").Value = '2.959'
$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("D34").Value = '4.480'
$ws.Range("E34").Value = '  -1.85%  '

$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  +2.60%  '

$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("E37").Value = '  +0.73%  '

$ws.Range("D38").Value = '0.05413'
$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("D39").Value = '0.01922'
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("D40").Value = '2.811'
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D41").Value = '0.1654'
$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("D42").Value = '0.5050'
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").Value = '6.526'
$ws.Range("E43").Value = '  -4.87%  '

$ws.Range("D44").Value = '8.349'
$ws.Range("E44").Value = '  +0.98%  '

$ws.Range("D45").Value = '0.06539'
$ws.Range("E45").Value = '  -0.94%  '

$ws.Range("D46").Value = '105.94'
$ws.Range("E46").Value = '  +1.41%  '

$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("D49").Value = '0.4624'
$ws.Range("E49").Value = '  -1.58%  '

$ws.Range("D50").Value = '1.629'
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").Value = '64.19'
$ws.Range("E51").Value = '  -0.40%  '
